$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (styles, fills, borders, row height, etc.) from the
# existing header/body rows onto the new rows being appended, then fill in
# the values/text for the new "Sprint" block (rows 7-11).

$ws.Range("A1:D1").Copy()
$ws.Range("A7:D7").PasteSpecial(-4122)

$ws.Range("A2:D2").Copy()
$ws.Range("A8:D8").PasteSpecial(-4122)

$ws.Range("A3:D3").Copy()
$ws.Range("A9:D9").PasteSpecial(-4122)

$ws.Range("A4:D4").Copy()
$ws.Range("A10:D10").PasteSpecial(-4122)

$ws.Range("A5:D5").Copy()
$ws.Range("A11:D11").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Row heights for the new rows ---
$ws.Rows.Item(7).RowHeight = 82.5
$ws.Rows.Item(8).RowHeight = 82.5
$ws.Rows.Item(9).RowHeight = 82.5
$ws.Rows.Item(10).RowHeight = 82.5
$ws.Rows.Item(11).RowHeight = 82.5

# --- Row 7: repeat header row ---
$ws.Range("A7").Value = 44244
$ws.Range("B7").Value = "1- Whats`n been done"
$ws.Range("C7").Value = "2- what `nwill we do"
$ws.Range("D7").Value = "3- difficulties`n encountered"

# --- Row 8: Anibal ---
$ws.Range("A8").Value = "Anibal"
$ws.Range("B8").Value = "UC2 - PLSQL"
$ws.Range("C8").Value = " Implementação Novas Classes TDD"
$ws.Range("D8").Value = "aplicacao de alguns dos conceitos BD"

# --- Row 9: Julio ---
$ws.Range("A9").Value = "Julio"
$ws.Range("B9").Value = "UC4, UC6 - PLSQL"
$ws.Range("C9").Value = "Finalização PLSQL + ligação modelo com base de dados"
$ws.Range("D9").Value = "aplicacao de alguns dos conceitos BD"

# --- Row 10: Pedro ---
$ws.Range("A10").Value = "Pedro"
$ws.Range("B10").Value = "UC1, UC3 - PLSQL`nUC1 - Ligação Modelo com base de dados"
$ws.Range("C10").Value = "Finalização PLSQL + ligação modelo com base de dados"
$ws.Range("D10").Value = "aplicacao de alguns dos conceitos BD"

# --- Row 11: Vitor ---
$ws.Range("A11").Value = "Vitor"
$ws.Range("B11").Value = "UC5 - PLSQL`nInserts Gerais"
$ws.Range("C11").Value = "Implementação Novas Classes TDD"
$ws.Range("D11").Value = "aplicacao de alguns dos conceitos BD"

# --- View / selection updates ---
$ws.Application.ActiveWindow.Zoom = 85
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("B8").Select()
